$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos list refresh.
# Values that look numeric are written with a leading apostrophe
# (quote-prefix) so Excel keeps them as text, matching the original
# inline-string cell type; the Style reset avoids leaving a stray
# quote-prefix style on the cell.

$ws.Range("D2").Value = '26.478.97'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.725.72'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D4").Value = '''0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''245.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.21%  '
$ws.Range("D6").Value = '''0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '''0.4810'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.20%  '
$ws.Range("D8").Value = '''0.2666'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.54%  '
$ws.Range("D9").Value = '''0.06183'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("D10").Value = '1.736.87'
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("D11").Value = '''0.07188'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '''15.57'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").Value = '''0.6106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("D14").Value = '''4.520'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.98%  '
$ws.Range("D15").Value = '''77.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = '''0.9996'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '26.498.36'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '''0.9997'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  +1.94%  '
$ws.Range("D20").Value = '''11.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").Value = '1.957.98'
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("D23").Value = '''8.791'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.71%  '
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").Value = '''136.91'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("D26").Value = '''15.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("D27").Value = '''1.778'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("D28").Value = '''1.393'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.56%  '
$ws.Range("D29").Value = '''107.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  -0.48%  '
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("D32").Value = '''3.687'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("D34").Value = '''0.9991'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = '''2.611'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.20%  '
$ws.Range("D36").Value = '''0.9948'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.90%  '
$ws.Range("D37").Value = '''0.6262'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '''0.9136'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.36%  '
$ws.Range("D39").Value = '''2.071'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.60%  '
$ws.Range("D40").Value = '''2.377'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  +0.08%  '
$ws.Range("D42").Value = '''103.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.18%  '
$ws.Range("D43").Value = '''0.01505'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("D44").Value = '''5.640'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.21%  '
$ws.Range("D45").Value = '''0.3861'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("D46").Value = '''6.974'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.06%  '
$ws.Range("D47").Value = '''0.1180'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("D48").Value = '''0.05351'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.49%  '
$ws.Range("D49").Value = '''30.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").Value = '''7.770'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("D51").Value = '''1.251'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.85%  '
